# Marksheet fix: handle float input without breaking stuff.
# Consolidates answers that had spilled into the D/E and G/H "extra"
# column-pairs back into the single A/B pair, recomputes the summary
# row, and drops the now-unused G:H columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Summary rows (10-12) ---------------------------------------------
# Row 10 "No." -> Right / Wrong / Not Attempt / Max
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B10").Value = 17
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 28

# Row 11 "Marking"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 "Total"
$ws.Range("B12").Value = 68
$ws.Range("E12").Value = "68/112"

# --- Answer grid: pull the D/E (rows 16-18) and A/B (rows 19+) pairs --
# back in line, painting the "Student Ans" cell with the correctStyle
# (green) format whenever it matches the "Correct Ans" cell, using an
# already-styled cell as the format source so no new style records are
# created.
$ws.Range("B10").Copy() | Out-Null

$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = "Option A"

$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Value = "Option C"

$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("D18").Value = "Option D"

$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Value = "Option C"

$ws.Range("A22").PasteSpecial(-4122) | Out-Null
$ws.Range("A22").Value = "Option D"

$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value = "Option D"

$ws.Range("A25").PasteSpecial(-4122) | Out-Null
$ws.Range("A25").Value = "Option A"

$ws.Range("A27").PasteSpecial(-4122) | Out-Null
$ws.Range("A27").Value = "Option A"

$ws.Range("A28").PasteSpecial(-4122) | Out-Null
$ws.Range("A28").Value = "Option D"

$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("A30").Value = "Option B"

$ws.Range("A31").PasteSpecial(-4122) | Out-Null
$ws.Range("A31").Value = "Option D"

$ws.Range("A32").PasteSpecial(-4122) | Out-Null
$ws.Range("A32").Value = "Option C"

$ws.Range("A33").PasteSpecial(-4122) | Out-Null
$ws.Range("A33").Value = "Option D"

$ws.Range("A36").PasteSpecial(-4122) | Out-Null
$ws.Range("A36").Value = "Option A"

$ws.Range("A38").PasteSpecial(-4122) | Out-Null
$ws.Range("A38").Value = "Option A"

$ws.Range("A39").PasteSpecial(-4122) | Out-Null
$ws.Range("A39").Value = "Option D"

$ws.Range("A40").PasteSpecial(-4122) | Out-Null
$ws.Range("A40").Value = "Option D"

$excel.CutCopyMode = 0

# The duplicate "Student Ans"/"Correct Ans" pair that used to live in
# D:E for rows 19 and below is no longer needed now that the real
# answers live in A:B.
$ws.Range("D19:E40").Clear()

# --- Drop the now-empty third answer-pair (G:H) ------------------------
$ws.Columns("G:H").Delete()
